# Adds the "Nome do Arquivo" column (H) and rewrites "Valor NF" (G) values
# from numeric amounts into pt-BR formatted text strings (e.g. 1.240,00).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = 'Nome do Arquivo'

$rows = @(
    @(2, '250,00', '2025.05.26_PLUXEE BENEF_4992069.pdf'),
    @(3, '2.000,00', 'COMPRA CREDITO VR FLEX _35826113- NF 04981167.pdf'),
    @(4, '1.240,00', 'L A SANTANA NFSe 210.pdf'),
    @(5, '100,00', 'L A SANTANA NFSe 211.pdf'),
    @(6, '21.417,20', 'NF 13788 - CASA BAHIA COMERCIAL LTDA.pdf'),
    @(7, '691,83', 'NF 158707.pdf'),
    @(8, '1.750,00', 'NF 1634.pdf'),
    @(9, '12.753,45', 'NF 303910.pdf'),
    @(10, '12.371,95', 'NF 303911.pdf'),
    @(11, '5.000,00', 'NF 7060.pdf'),
    @(12, '645,95', 'NF 86746 2A SISTEMA AMBIENTAL LTDA.pdf'),
    @(13, '12.371,95', 'NF-e - Nota Fiscal Eletrônica de Serviços - São Paulo-303911-boleto.pdf'),
    @(14, '13.564,38', 'NF-e - Nota Fiscal Eletrônica de Serviços - São Paulo-AAEIFS.pdf'),
    @(15, '1.197,70', 'NF-e - Nota Fiscal Eletrônica de Serviços - São Paulo-AAKCID.pdf'),
    @(16, '691,83', 'NF-e - Nota Fiscal Eletrônica de Serviços - São Paulo.pdf'),
    @(17, '9.561,45', 'NF-e 158706.pdf'),
    @(18, '6.630,98', 'NF820 CASAS BAHIA R$6.630,98.pdf'),
    @(19, '2.250,00', 'NF823 CASAS BAHIA R$2.250,00.pdf'),
    @(20, '600,02', 'NFSe22340_40712883000187.pdf'),
    @(21, '5.000,00', 'NFSe3765_02878522000116.pdf'),
    @(22, '27.533,95', 'NFSe5661_03558771000197.pdf'),
    @(23, '1.161,94', 'NFSe_411_ICONREALTY_PLFALESIM.pdf'),
    @(24, '300,00', 'Nota 2422.pdf'),
    @(25, '5.600,00', 'nota_329557 (NFSe).pdf'),
    @(26, '5.600,00', 'nota_336834 (NFSe).pdf'),
    @(27, '15.000,00', 'N_2025.05.26_000_7507_4500012192.pdf'),
    @(28, '254,00', 'N_2025.06.03_000_10792_4500012306.pdf'),
    @(29, '234,00', 'N_2025.06.03_000_10802_4500012330.pdf'),
    @(30, '201,25', 'N_2025.06.03_009_40839_4500011946.pdf'),
    @(31, '250,00', 'RETHA ITAQUERA - NF 25269 - MAIO.pdf')
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $valorNf = $r[1]
    $nomeArquivo = $r[2]
    $ws.Range("G$rowNum").Value = $valorNf
    $ws.Range("H$rowNum").Value = $nomeArquivo
}
